$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 8
$ws.Range("B5").Value = 20
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 100

$ws.Range("B4").Select()
